$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.182.57"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3
$ws.Range("D3").Value = "3.773.74"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'630.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.84%  "

# Row 6
$ws.Range("D6").Value = "'166.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.95%  "

# Row 7
$ws.Range("D7").Value = "3.772.51"
$ws.Range("E7").Value = "  -0.84%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.67%  "

# Row 10
$ws.Range("E10").Value = "  -0.08%  "

# Row 11
$ws.Range("D11").Value = "'0.460"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.89%  "

# Row 12
$ws.Range("D12").Value = "'6.77"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.58%  "

# Row 13
$ws.Range("E13").Value = "  -1.58%  "

# Row 14
$ws.Range("D14").Value = "'34.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.68%  "

# Row 15
$ws.Range("D15").Value = "4.407.04"
$ws.Range("E15").Value = "  -0.72%  "

# Row 16
$ws.Range("D16").Value = "3.778.32"
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("D17").Value = "69.165.63"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18
$ws.Range("D18").Value = "'17.57"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.16%  "

# Row 19
$ws.Range("E19").Value = "  -0.87%  "

# Row 20
$ws.Range("D20").Value = "'7.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.94%  "

# Row 21
$ws.Range("D21").Value = "'464.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.38%  "

# Row 22
$ws.Range("D22").Value = "'9.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.45%  "

# Row 23
$ws.Range("D23").Value = "'0.708"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.33%  "

# Row 24
$ws.Range("D24").Value = "'83.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.80%  "

# Row 25
$ws.Range("D25").Value = "'0.0000145"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.67%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'2.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.01%  "

# Row 28
$ws.Range("D28").Value = "'10.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.44%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").Value = "3.922.10"
$ws.Range("E30").Value = "  -0.62%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.19%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.68%  "

# Row 33
$ws.Range("D33").Value = "'7.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.94%  "

# Row 34
$ws.Range("D34").Value = "'28.61"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.03%  "

# Row 35
$ws.Range("D35").Value = "'0.176"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +16.69%  "

# Row 36
$ws.Range("E36").Value = "  +0.19%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'8.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.25%  "

# Row 38
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.723.87"
$ws.Range("E38").Value = "  -0.75%  "

# Row 39
$ws.Range("E39").Value = "  +0.88%  "

# Row 40
$ws.Range("D40").Value = "'3.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.97%  "

# Row 41
$ws.Range("D41").Value = "'5.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.34%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.18%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.961"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.09%  "

# Row 44
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "'158.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.48%  "

# Row 46
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.55%  "

# Row 47
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'43.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("D48").Value = "'1.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.07%  "

# Row 49
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.296"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.47%  "

# Row 50
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'46.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.62%  "

# Row 51
$ws.Range("D51").Value = "'8.37"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.28%  "
